# Adventurers-Inc Names.xlsx: separate movement from generation, modify click routine.
# - human_male: "Eiikichi" renamed to "Rolph"; "Moredecaï" renamed/fixed to "Mordecaï"
# - elf_female: "Yalandhadra" renamed/fixed to "Yalandhra"
# - cursor/selection moved on a few sheets as a side effect of editing

$wb = $excel.ActiveWorkbook

# --- human_male: row 9 (Eiikichi -> Rolph), row 20 (Moredecaï -> Mordecaï) ---
$wsHumanMale = $wb.Worksheets.Item("human_male")
$wsHumanMale.Activate() | Out-Null
$wsHumanMale.Range("A9").Value = "Rolph"
$wsHumanMale.Range("A20").Value = "Mordecaï"
$wsHumanMale.Range("D25").Select() | Out-Null

# --- elf_female: row 10 (Yalandhadra -> Yalandhra) ---
$wsElfFemale = $wb.Worksheets.Item("elf_female")
$wsElfFemale.Activate() | Out-Null
$wsElfFemale.Range("A10").Value = "Yalandhra"
$wsElfFemale.Range("I11").Select() | Out-Null

# --- elf_neutral: move the selection/cursor (stays the active sheet) ---
$wsElfNeutral = $wb.Worksheets.Item("elf_neutral")
$wsElfNeutral.Activate() | Out-Null
$wsElfNeutral.Range("E9").Select() | Out-Null
